# Updates the cryptocurrency price/volume snapshot on Sheet1 to the
# refreshed values published by the scraper run.
#
# Rows 44 and 45 also swap their Coin/Link contents (EnergySwap and
# FirstDigitalUSD traded ranking positions), so those two rows update
# columns B-E instead of just D/E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "64.273.59"; Numeric = $false },
    @{ Cell = "E2"; Value = "  +0.93%  "; Numeric = $false },
    @{ Cell = "D3"; Value = "3.484.38"; Numeric = $false },
    @{ Cell = "E3"; Value = "  +0.59%  "; Numeric = $false },
    @{ Cell = "E4"; Value = "  +0.04%  "; Numeric = $false },
    @{ Cell = "D5"; Value = "586.04"; Numeric = $true },
    @{ Cell = "E5"; Value = "  +0.80%  "; Numeric = $false },
    @{ Cell = "D6"; Value = "133.83"; Numeric = $true },
    @{ Cell = "E6"; Value = "  +2.27%  "; Numeric = $false },
    @{ Cell = "D7"; Value = "3.484.49"; Numeric = $false },
    @{ Cell = "E7"; Value = "  +0.56%  "; Numeric = $false },
    @{ Cell = "E8"; Value = "  +0.02%  "; Numeric = $false },
    @{ Cell = "D9"; Value = "0.484"; Numeric = $true },
    @{ Cell = "E9"; Value = "  -0.48%  "; Numeric = $false },
    @{ Cell = "E10"; Value = "  +0.55%  "; Numeric = $false },
    @{ Cell = "D11"; Value = "7.20"; Numeric = $true },
    @{ Cell = "E11"; Value = "  +1.76%  "; Numeric = $false },
    @{ Cell = "E12"; Value = "  -2.12%  "; Numeric = $false },
    @{ Cell = "D13"; Value = "4.081.15"; Numeric = $false },
    @{ Cell = "E13"; Value = "  +0.82%  "; Numeric = $false },
    @{ Cell = "E14"; Value = "  +2.22%  "; Numeric = $false },
    @{ Cell = "D15"; Value = "0.0000179"; Numeric = $true },
    @{ Cell = "E15"; Value = "  +1.64%  "; Numeric = $false },
    @{ Cell = "D16"; Value = "3.486.84"; Numeric = $false },
    @{ Cell = "E16"; Value = "  +0.70%  "; Numeric = $false },
    @{ Cell = "D17"; Value = "64.322.38"; Numeric = $false },
    @{ Cell = "E17"; Value = "  +0.97%  "; Numeric = $false },
    @{ Cell = "D18"; Value = "25.25"; Numeric = $true },
    @{ Cell = "E18"; Value = "  -8.22%  "; Numeric = $false },
    @{ Cell = "D19"; Value = "9.99"; Numeric = $true },
    @{ Cell = "E19"; Value = "  +1.33%  "; Numeric = $false },
    @{ Cell = "D20"; Value = "5.68"; Numeric = $true },
    @{ Cell = "E20"; Value = "  +0.87%  "; Numeric = $false },
    @{ Cell = "D21"; Value = "13.65"; Numeric = $true },
    @{ Cell = "E21"; Value = "  -3.91%  "; Numeric = $false },
    @{ Cell = "D22"; Value = "383.52"; Numeric = $true },
    @{ Cell = "E22"; Value = "  -1.51%  "; Numeric = $false },
    @{ Cell = "D23"; Value = "0.565"; Numeric = $true },
    @{ Cell = "E23"; Value = "  -1.24%  "; Numeric = $false },
    @{ Cell = "D24"; Value = "3.624.44"; Numeric = $false },
    @{ Cell = "E24"; Value = "  +0.60%  "; Numeric = $false },
    @{ Cell = "D25"; Value = "74.07"; Numeric = $true },
    @{ Cell = "E25"; Value = "  +1.70%  "; Numeric = $false },
    @{ Cell = "E26"; Value = "  +0.04%  "; Numeric = $false },
    @{ Cell = "D27"; Value = "5.69"; Numeric = $true },
    @{ Cell = "E27"; Value = "  -0.65%  "; Numeric = $false },
    @{ Cell = "E28"; Value = "  +4.89%  "; Numeric = $false },
    @{ Cell = "D29"; Value = "1.54"; Numeric = $true },
    @{ Cell = "E29"; Value = "  +0.42%  "; Numeric = $false },
    @{ Cell = "D30"; Value = "1.00"; Numeric = $true },
    @{ Cell = "E30"; Value = "  +0.29%  "; Numeric = $false },
    @{ Cell = "D31"; Value = "7.43"; Numeric = $true },
    @{ Cell = "E31"; Value = "  +0.60%  "; Numeric = $false },
    @{ Cell = "E32"; Value = "  +0.02%  "; Numeric = $false },
    @{ Cell = "D33"; Value = "8.19"; Numeric = $true },
    @{ Cell = "E33"; Value = "  +0.90%  "; Numeric = $false },
    @{ Cell = "D34"; Value = "3.505.53"; Numeric = $false },
    @{ Cell = "E34"; Value = "  +1.26%  "; Numeric = $false },
    @{ Cell = "E36"; Value = "  +2.48%  "; Numeric = $false },
    @{ Cell = "E37"; Value = "  -0.91%  "; Numeric = $false },
    @{ Cell = "D38"; Value = "5.28"; Numeric = $true },
    @{ Cell = "E38"; Value = "  +0.89%  "; Numeric = $false },
    @{ Cell = "E39"; Value = "  -1.22%  "; Numeric = $false },
    @{ Cell = "E40"; Value = "  -1.28%  "; Numeric = $false },
    @{ Cell = "D41"; Value = "162.39"; Numeric = $true },
    @{ Cell = "E41"; Value = "  -4.03%  "; Numeric = $false },
    @{ Cell = "E42"; Value = "  -2.81%  "; Numeric = $false },
    @{ Cell = "D43"; Value = "0.802"; Numeric = $true },
    @{ Cell = "E43"; Value = "  -0.47%  "; Numeric = $false },
    @{ Cell = "B44"; Value = "FirstDigitalUSD"; Numeric = $false },
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; Numeric = $false },
    @{ Cell = "D44"; Value = "1.00"; Numeric = $true },
    @{ Cell = "E44"; Value = "  +0.14%  "; Numeric = $false },
    @{ Cell = "B45"; Value = "EnergySwap"; Numeric = $false },
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; Numeric = $false },
    @{ Cell = "D45"; Value = "25.40"; Numeric = $true },
    @{ Cell = "E45"; Value = "  -0.34%  "; Numeric = $false },
    @{ Cell = "D46"; Value = "41.71"; Numeric = $true },
    @{ Cell = "E46"; Value = "  +0.39%  "; Numeric = $false },
    @{ Cell = "E47"; Value = "  +1.50%  "; Numeric = $false },
    @{ Cell = "D49"; Value = "1.64"; Numeric = $true },
    @{ Cell = "E49"; Value = "  +1.50%  "; Numeric = $false },
    @{ Cell = "D50"; Value = "2.461.26"; Numeric = $false },
    @{ Cell = "E50"; Value = "  +2.41%  "; Numeric = $false },
    @{ Cell = "E51"; Value = "  -1.50%  "; Numeric = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Numeric) {
        # Force text storage so Excel doesn't reinterpret values such as
        # "586.04" or "1.00" as numbers (the source data are plain strings).
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}
